$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.247.89"
$ws.Range("E2").Value = "  +2.57%  "
$ws.Range("D3").Value = "'2.975.85"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'564.72"
$ws.Range("E5").Value = "  +1.95%  "
$ws.Range("D6").Value = "'137.25"
$ws.Range("E6").Value = "  +2.97%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.518"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").Value = "'2.967.65"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("D11").Value = "'5.33"
$ws.Range("E11").Value = "  +10.42%  "
$ws.Range("D12").Value = "'0.449"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "'0.0000228"
$ws.Range("E13").Value = "  +3.40%  "
$ws.Range("D14").Value = "'33.51"
$ws.Range("E14").Value = "  +2.30%  "
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "'3.472.68"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").Value = "'7.04"
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("D18").Value = "'2.978.89"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("D19").Value = "'59.267.72"
$ws.Range("E19").Value = "  +2.66%  "
$ws.Range("D20").Value = "'435.20"
$ws.Range("E20").Value = "  +4.34%  "
$ws.Range("D21").Value = "'13.54"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D22").Value = "'0.719"
$ws.Range("E22").Value = "  +2.32%  "
$ws.Range("D23").Value = "'6.98"
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("D24").Value = "'13.12"
$ws.Range("E24").Value = "  -2.00%  "
$ws.Range("D25").Value = "'79.78"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D27").Value = "'2.21"
$ws.Range("E27").Value = "  +9.06%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").Value = "'2.54"
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("D30").Value = "'7.69"
$ws.Range("E30").Value = "  +2.36%  "
$ws.Range("D31").Value = "'6.20"
$ws.Range("E31").Value = "  +4.52%  "
$ws.Range("D32").Value = "'25.63"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("E33").Value = "  +7.58%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.87"
$ws.Range("E34").Value = "  +2.80%  "
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "'0.0₃0763"
$ws.Range("E35").Value = "  +8.62%  "
$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").Value = "'0.979"
$ws.Range("E36").Value = "  +2.61%  "
$ws.Range("D37").Value = "'2.06"
$ws.Range("E37").Value = "  +0.62%  "
$ws.Range("D38").Value = "'48.49"
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("D39").Value = "'8.68"
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("E40").Value = "  +2.86%  "
$ws.Range("D41").Value = "'397.36"
$ws.Range("E41").Value = "  +3.13%  "
$ws.Range("D42").Value = "'0.0350"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").Value = "'2.724.67"
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").Value = "'0.104"
$ws.Range("E44").Value = "  -2.62%  "
$ws.Range("D45").Value = "'0.249"
$ws.Range("E45").Value = "  +5.06%  "
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'122.19"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").Value = "'33.89"
$ws.Range("E48").Value = "  +16.30%  "
$ws.Range("D49").Value = "'0.109"
$ws.Range("E49").Value = "  +1.58%  "
$ws.Range("D50").Value = "'1.98"
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("D51").Value = "'23.10"
$ws.Range("E51").Value = "  +1.48%  "
